# Apply the template-cleanup edit:
#   - Remove the second slide ("Template editing instructions and feedback"),
#     which also frees up its slide layout.
#   - Remove the now-unused "1_Blank" custom layout (Date/Footer/SlideNumber
#     placeholders) that only the deleted slide referenced.

$p = $ppt.ActivePresentation

# --- Remove slide 2 ("Template editing instructions") -----------------
# Slide 1 is the main dashboard (layout "Chart dashboard"); slide 2 is the
# standalone instructions slide built on the "1_Blank" layout.
if ($p.Slides.Count -ge 2) {
    $p.Slides.Item(2).Delete()
}

# --- Remove the now-orphaned "1_Blank" custom layout -------------------
$master = $p.SlideMaster
for ($i = $master.CustomLayouts.Count; $i -ge 1; $i--) {
    $layout = $master.CustomLayouts.Item($i)

    $inUse = $false
    for ($j = 1; $j -le $p.Slides.Count; $j++) {
        if ($p.Slides.Item($j).CustomLayout.Index -eq $layout.Index) {
            $inUse = $true
        }
    }

    if (-not $inUse -and $layout.Shapes.Count -gt 0) {
        $hasDatePh = $false
        for ($k = 1; $k -le $layout.Shapes.Count; $k++) {
            if ($layout.Shapes.Item($k).Name -eq "Date Placeholder 1") {
                $hasDatePh = $true
            }
        }
        if ($hasDatePh) {
            $layout.Delete()
        }
    }
}
